{"js": "// Office.js (Word JavaScript API) implementation of the diff:\n//  1. The paragraph that only held the \"_GoBack\" bookmark (right after the\n//     title) gains the \"\u0412\u044b\u043f\u043e\u043b\u043d\u0438\u043b\u0430: ...\" text (moved up from the paragraph\n//     below it), is followed by a new empty paragraph, and then a new\n//     paragraph \"\u0421\u0441\u044b\u043b\u043a\u0430 \u043d\u0430 \u0440\u0435\u043f\u043e\u0437\u0438\u0442\u043e\u0440\u0438\u0439: https://...\" which now carries the\n//     \"_GoBack\" bookmark at its end.\n//  2. The old \"\u0412\u044b\u043f\u043e\u043b\u043d\u0438\u043b\u0430: ...\" paragraph becomes empty.\n//  3. The \"#include <stdio.h>\" paragraph has its first three runs\n//     (\"#\", \"include\", \" <\") merged into a single run and gains\n//     w:lang=\"en-US\" on the paragraph mark and on every run.\n\nconst runPr = '<w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr>';\n\n// Package helper: wraps a <w:body> fragment into a full OOXML \"flat\" part\n// suitable for Paragraph.insertOoxml().\nfunction pkg(bodyXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\n// --- Step 1: locate the bookmark-only paragraph right after the title. ---\nlet paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst bookmarkParaIndex = 1; // paragraph right after the title\nconst oldAuthorParaIndex = 2; // paragraph holding \"\u0412\u044b\u043f\u043e\u043b\u043d\u0438\u043b\u0430: ...\"\n\nconst bookmarkPara = paragraphs.items[bookmarkParaIndex];\n\nconst newIntroBlock =\n  '<w:p><w:pPr>' + runPr + '</w:pPr>' +\n  '<w:r>' + runPr + '<w:t xml:space=\"preserve\">\u0412\u044b\u043f\u043e\u043b\u043d\u0438\u043b\u0430: \u0441\u0442\u0443\u0434\u0435\u043d\u0442\u043a\u0430 \u0433\u0440\u0443\u043f\u043f\u044b 213-723 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r>' + runPr + '<w:t>\u041a\u0438\u044f\u0447\u0435\u043d\u043a\u043e</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r>' + runPr + '<w:t xml:space=\"preserve\"> \u0412\u0435\u0440\u0430 \u0410\u043d\u0434\u0440\u0435\u0435\u0432\u043d\u0430</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr>' + runPr + '</w:pPr></w:p>' +\n  '<w:p><w:pPr>' + runPr + '</w:pPr>' +\n  '<w:r>' + runPr + '<w:t xml:space=\"preserve\">\u0421\u0441\u044b\u043b\u043a\u0430 \u043d\u0430 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r>' + runPr + '<w:t>\u0440\u0435\u043f\u043e\u0437\u0438\u0442\u043e\u0440\u0438\u0439</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r>' + runPr + '<w:t xml:space=\"preserve\">: </w:t></w:r>' +\n  '<w:r>' + runPr + '<w:t>https://github.com/privetverok/polytech-introduction-to-programing</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>';\n\nbookmarkPara.insertOoxml(pkg(newIntroBlock), Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 2: the content that used to live right under the bookmark\n// paragraph (\"\u0412\u044b\u043f\u043e\u043b\u043d\u0438\u043b\u0430: ...\") moved up, so that paragraph becomes empty.\n// Re-fetch paragraphs since the collection shifted (1 paragraph -> 3, a net\n// +2 growth), so the old author paragraph is now 2 positions further down.\n// (The new intro paragraph has identical text, so a text search would find\n// the wrong match first \u2014 the index shift is deterministic, use it.)\nparagraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst oldAuthorPara = paragraphs.items[oldAuthorParaIndex + 2];\n\noldAuthorPara.insertOoxml(pkg('<w:p><w:pPr>' + runPr + '</w:pPr></w:p>'), Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 3: merge the \"#include <stdio.h>\" runs and tag them en-US. ---\n// This exact line repeats several times later in the document (as sample\n// code), so it must be located positionally, not by a text search: the\n// original target was paragraph index 4, and the net +2 paragraphs added\n// above shift it to index 6.\nparagraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst includeParaIndex = 4 + 2;\nconst includePara = paragraphs.items[includeParaIndex];\n\nconst includeRpr =\n  '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:lang w:val=\"en-US\"/></w:rPr>';\n\nconst newIncludeParaXml =\n  '<w:p><w:pPr><w:spacing w:line=\"360\" w:lineRule=\"auto\"/>' +\n  '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n  '<w:r>' + includeRpr + '<w:t>#include &lt;</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r>' + includeRpr + '<w:t>stdio.h</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r>' + includeRpr + '<w:t>&gt;</w:t></w:r>' +\n  '</w:p>';\n\nincludePara.insertOoxml(pkg(newIncludeParaXml), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop implementation of the diff:\n#  1. The paragraph that only held the \"_GoBack\" bookmark (right after the\n#     title) gains the \"\u0412\u044b\u043f\u043e\u043b\u043d\u0438\u043b\u0430: ...\" text (moved up from the paragraph\n#     below it), is followed by a new empty paragraph, and then a new\n#     paragraph \"\u0421\u0441\u044b\u043b\u043a\u0430 \u043d\u0430 \u0440\u0435\u043f\u043e\u0437\u0438\u0442\u043e\u0440\u0438\u0439: https://...\" which now carries the\n#     \"_GoBack\" bookmark at its end.\n#  2. The old \"\u0412\u044b\u043f\u043e\u043b\u043d\u0438\u043b\u0430: ...\" paragraph becomes empty.\n#  3. The \"#include <stdio.h>\" paragraph has its first three runs\n#     (\"#\", \"include\", \" <\") merged into a single run and gains\n#     w:lang=\"en-US\" on the paragraph mark and on every run.\n\n$d = $word.ActiveDocument\n\n$pkgOpen = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$runPr = '<w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr>'\n\n# --- Step 1: paragraph right after the title only held the bookmark. ---\n$bookmarkParaIndex = 2   # Word COM Paragraphs is 1-based; title is #1\n$oldAuthorParaIndex = 3  # paragraph holding \"\u0412\u044b\u043f\u043e\u043b\u043d\u0438\u043b\u0430: ...\"\n\n$newIntroBlock = '<w:p><w:pPr>' + $runPr + '</w:pPr>' + `\n  '<w:r>' + $runPr + '<w:t xml:space=\"preserve\">\u0412\u044b\u043f\u043e\u043b\u043d\u0438\u043b\u0430: \u0441\u0442\u0443\u0434\u0435\u043d\u0442\u043a\u0430 \u0433\u0440\u0443\u043f\u043f\u044b 213-723 </w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellStart\"/>' + `\n  '<w:r>' + $runPr + '<w:t>\u041a\u0438\u044f\u0447\u0435\u043d\u043a\u043e</w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellEnd\"/>' + `\n  '<w:r>' + $runPr + '<w:t xml:space=\"preserve\"> \u0412\u0435\u0440\u0430 \u0410\u043d\u0434\u0440\u0435\u0435\u0432\u043d\u0430</w:t></w:r>' + `\n  '</w:p>' + `\n  '<w:p><w:pPr>' + $runPr + '</w:pPr></w:p>' + `\n  '<w:p><w:pPr>' + $runPr + '</w:pPr>' + `\n  '<w:r>' + $runPr + '<w:t xml:space=\"preserve\">\u0421\u0441\u044b\u043b\u043a\u0430 \u043d\u0430 </w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellStart\"/>' + `\n  '<w:r>' + $runPr + '<w:t>\u0440\u0435\u043f\u043e\u0437\u0438\u0442\u043e\u0440\u0438\u0439</w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellEnd\"/>' + `\n  '<w:r>' + $runPr + '<w:t xml:space=\"preserve\">: </w:t></w:r>' + `\n  '<w:r>' + $runPr + '<w:t>https://github.com/privetverok/polytech-introduction-to-programing</w:t></w:r>' + `\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' + `\n  '</w:p>'\n\n$bookmarkRange = $d.Paragraphs.Item($bookmarkParaIndex).Range\n$null = $bookmarkRange.InsertXML($pkgOpen + $newIntroBlock + $pkgClose)\n\n# --- Step 2: the content that used to live right under the bookmark\n# paragraph (\"\u0412\u044b\u043f\u043e\u043b\u043d\u0438\u043b\u0430: ...\") moved up, so that paragraph becomes empty.\n# The 1 paragraph we replaced above became 3, a net +2 growth, so the old\n# author paragraph is now 2 positions further down.\n$emptyParaXml = '<w:p><w:pPr>' + $runPr + '</w:pPr></w:p>'\n$oldAuthorRange = $d.Paragraphs.Item($oldAuthorParaIndex + 2).Range\n$null = $oldAuthorRange.InsertXML($pkgOpen + $emptyParaXml + $pkgClose)\n\n# --- Step 3: merge the \"#include <stdio.h>\" runs and tag them en-US. ---\n# This exact line repeats several times later in the document (as sample\n# code), so it must be located positionally: the original target was\n# paragraph #5 (1-based), and the net +2 paragraphs added above shift it to\n# paragraph #7.\n$includeParaIndex = 5 + 2\n$includeRpr = '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:lang w:val=\"en-US\"/></w:rPr>'\n\n$newIncludeParaXml = '<w:p><w:pPr><w:spacing w:line=\"360\" w:lineRule=\"auto\"/>' + `\n  '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' + `\n  '<w:r>' + $includeRpr + '<w:t>#include &lt;</w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellStart\"/>' + `\n  '<w:r>' + $includeRpr + '<w:t>stdio.h</w:t></w:r>' + `\n  '<w:proofErr w:type=\"spellEnd\"/>' + `\n  '<w:r>' + $includeRpr + '<w:t>&gt;</w:t></w:r>' + `\n  '</w:p>'\n\n$includeRange = $d.Paragraphs.Item($includeParaIndex).Range\n$null = $includeRange.InsertXML($pkgOpen + $newIncludeParaXml + $pkgClose)\n"}
